$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Add-RequiredMarker {
    param(
        [string]$CellAddr,
        [bool]$MarkerIncludesSpace
    )
    $cell = $ws.Range($CellAddr)
    $orig = $cell.Text
    $cell.Value = $orig + " (*)"
    # Bold the whole cell (base text inherits the cell's bold style,
    # the trailing "(*)" run gets overridden to red below).
    $cell.Font.Bold = $true
    $newText = $cell.Text
    $len = $newText.Length
    if ($MarkerIncludesSpace) {
        $marker = $cell.Characters($len - 3, 4)
    } else {
        $marker = $cell.Characters($len - 2, 3)
    }
    $marker.Font.Color = 255   # BGR-encoded RGB(255,0,0) == red
}

Add-RequiredMarker "A1" $false
Add-RequiredMarker "B1" $false
Add-RequiredMarker "C1" $false
Add-RequiredMarker "D1" $true
Add-RequiredMarker "E1" $true

# Columns C:E auto-size their width based on content (bestFit); re-fit them
# now that the header text is longer.
$ws.Columns("C:E").AutoFit() | Out-Null

# Update the saved selection, as recorded in the edited workbook.
$ws.Range("C9").Select() | Out-Null
